$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Merge the split runs in step "2" of the "Flujo de Sucesos" table into a
#    single run (text content is unchanged, only run boundaries collapse).
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    "El usuario ingresa # (numeral/almohadilla) seguido de la palabra buscada.          {cd ",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "El usuario ingresa # (numeral/almohadilla) seguido de la palabra buscada.          {cd ",
    2) | Out-Null

# ---------------------------------------------------------------------------
# 2) Merge the split runs in step "3" ("El sistema filtra los mensajes por la
#    etiqueta ingresada") into a single run.
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    "El sistema filtra los mensajes por la etiqueta ingresada",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "El sistema filtra los mensajes por la etiqueta ingresada",
    2) | Out-Null

# ---------------------------------------------------------------------------
# 3) Merge the split runs in step "4" ("El usuario visualiza los mensajes que
#    poseen dicha etiqueta.") into a single run.
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    "El usuario visualiza los mensajes que poseen dicha etiqueta.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "El usuario visualiza los mensajes que poseen dicha etiqueta.",
    2) | Out-Null

# ---------------------------------------------------------------------------
# 4) Collapse the two trailing empty paragraphs (after the flow table, before
#    the "Subflujos" table) into one, and replace its content with a manual
#    page break run (dropping the explicit spacing override on the
#    paragraph). NOTE: this must run *before* any table-row deletion below —
#    the document-wide Paragraphs collection gets confused once a table row
#    has been removed earlier in the body.
# ---------------------------------------------------------------------------
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -eq [char]13 -and $p.Next().Range.Text -eq [char]13) {
        $nextp = $p.Next()
        if ($nextp.Next().Range.Text -like "Subflujos*") {
            $mergeRange = $d.Range($p.Range.Start, $nextp.Range.End)
            $mergeRange.Delete()

            $target = $d.Paragraphs.Item($i)
            $fullPara = $d.Range($target.Range.Start, $target.Range.End)

            $xml = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
<w:p>
  <w:pPr>
    <w:rPr>
      <w:rFonts w:ascii="Times New Roman" w:eastAsia="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>
      <w:sz w:val="24"/>
      <w:szCs w:val="24"/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Times New Roman" w:eastAsia="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>
      <w:sz w:val="24"/>
      <w:szCs w:val="24"/>
    </w:rPr>
    <w:br w:type="page"/>
  </w:r>
</w:p>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
'@

            $fullPara.InsertXML($xml) | Out-Null
            break
        }
    }
}

# ---------------------------------------------------------------------------
# 5) Remove the final "5 / Fin de caso de uso." row from the flow table.
# ---------------------------------------------------------------------------
for ($i = 1; $i -le $d.Tables.Count; $i++) {
    $tbl = $d.Tables.Item($i)
    $lastRow = $tbl.Rows.Item($tbl.Rows.Count)
    if ($lastRow.Cells.Item($lastRow.Cells.Count).Range.Text -like "Fin de caso de uso.*") {
        $lastRow.Delete()
        break
    }
}
